$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44473
$ws.Range("M3").Value = 120

$ws.Range("D4").Value = 44438
$ws.Range("M4").Value = 60

$ws.Range("D5").Value = 44432
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 1300
$ws.Range("O5").Value = 1300
$ws.Range("P5").Value = 1300
$ws.Range("S5").Value = 1300

$ws.Range("D6").Value = 44476
$ws.Range("M6").Value = 80

$ws.Range("D7").Value = 44431
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 1300
$ws.Range("O7").Value = 1300
$ws.Range("P7").Value = 1300
$ws.Range("S7").Value = 1300

$ws.Range("D8").Value = 44424
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 1200
$ws.Range("O8").Value = 1200
$ws.Range("P8").Value = 1200
$ws.Range("S8").Value = 1200

$ws.Range("D9").Value = 44418
$ws.Range("M9").Value = 40

$ws.Range("D10").Value = 44405
$ws.Range("M10").Value = 50

$ws.Range("D11").Value = 44357
$ws.Range("M11").Value = 35
$ws.Range("N11").Value = 1000
$ws.Range("O11").Value = 1000
$ws.Range("P11").Value = 1000
$ws.Range("S11").Value = 1000

$ws.Range("D12").Value = 44343
$ws.Range("M12").Value = 60

$ws.Range("D13").Value = 44417
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 1200
$ws.Range("O13").Value = 1200
$ws.Range("P13").Value = 1200
$ws.Range("S13").Value = 1200
